$wb = $excel.ActiveWorkbook

# Both the "展览" (exhibitions) sheet and the "全部类型" (all types) sheet
# contain the same four data rows (rows 2-5) whose "想去人数" (F column)
# counts need to be reset to 0.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2:F5").Value = 0
}
